# Weekly update: a new price observation was collected for
# "Terminal La Palmera de La Serena - Espinaca" and inserted into the
# dataset at row 53 (it is the new most-recent-fetched record while the
# sheet stays ordered as previously), pushing every following record
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 53; everything from the old row 53
# onward (old rows 53-171) shifts down to become rows 54-172.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A53").Value = 8
$ws.Range("B53").Value = "Terminal La Palmera de La Serena"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44498
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = 100112012
$ws.Range("G53").Value = "Espinaca"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 3200
$ws.Range("K53").Value = 400
$ws.Range("L53").Value = 500
$ws.Range("M53").Value = 450
$ws.Range("N53").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O53").Value = "Provincia del Elquí"
$ws.Range("P53").Value = 900
$ws.Range("Q53").Value = 0.5
$ws.Range("R53").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D53").NumberFormat = $ws.Range("D54").NumberFormat
